# Auto-generated Word COM-interop script
# Adds '2 Corinthians Ch5' verse list (and the Ch6 heading stub)
# after the existing '2 Corinthians Ch5' title paragraph.

$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

# --- locate the existing chapter-5 title paragraph (last body paragraph) ---
$titlePara = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "2 Corinthians Ch5") {
        $titlePara = $d.Paragraphs.Item($i)
        break
    }
}
if ($null -eq $titlePara) {
    throw "Could not find the '2 Corinthians Ch5' title paragraph"
}
$anchorIndex = $titlePara.Index

# --- bootstrap paragraph used only to mint a fresh numbered-list definition ---
$d.Paragraphs.Item($anchorIndex).Range.InsertParagraphAfter()
$bootstrapIndex = $anchorIndex + 1
$bootstrap = $d.Paragraphs.Item($bootstrapIndex)
$bootstrap.Range.Text = "x"
$bootstrap.Range.ListFormat.ApplyNumberDefault()

# Recover the numId that got minted for the bootstrap paragraph.
$mintedXml = $d.Paragraphs.Item($bootstrapIndex).Range.WordOpenXML
if ($mintedXml -match 'w:numId w:val="(\d+)"') {
    $numId = [int]$Matches[1]
} else {
    throw "Could not recover minted numId"
}

# Normalize the minted abstract numbering definition to the usual Word
# decimal / lowerLetter / lowerRoman repeating pattern (cosmetic parity
# with what Word mints for a fresh outline-numbered list).
$lt = $d.Paragraphs.Item($bootstrapIndex).Range.ListFormat.ListTemplate
for ($lvl = 1; $lvl -le 9; $lvl++) {
    $m = ($lvl - 1) % 3
    $level = $lt.ListLevels.Item($lvl)
    if ($m -eq 1) {
        $level.NumberStyle = 4
        $level.Alignment = 0
    } elseif ($m -eq 2) {
        $level.NumberStyle = 2
        $level.Alignment = 2
    } else {
        $level.NumberStyle = 0
        $level.Alignment = 0
    }
}

# --- now lay down the real paragraphs, each stamped with exact OOXML ---
$prevIndex = $bootstrapIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">For we know that if our earthly house of this tabernacle were dissolved, we have a building of God, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>an</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> house not made with hands, eternal in the heavens.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>For in this we groan, earnestly desiring to be clothed upon with our house which is from heaven:</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">If </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>so</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> be that being clothed we shall not be found naked.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>For we that are in this tabernacle do groan, being burdened: not for that we would be unclothed, but clothed upon, that mortality might be swallowed up of life.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Now he that hath wrought us for the selfsame thing is God, who also hath given unto us the earnest of the Spirit.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Therefore</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> we are always confident, knowing that, whilst we are at home in the body, we are absent from the Lord:</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(For we walk by faith, not by sight:)</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>We are confident, I say, and willing rather to be absent from the body, and to be present with the Lord.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Wherefore we </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>labour</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, that, whether present or absent, we may be accepted of him.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">For we must all appear before the judgment seat of Christ; that </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>every one</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> may receive the things done in his body, according to that he hath done, whether it be good or bad.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Knowing therefore the terror of the Lord, we persuade men; but we are made manifest unto God; and I trust also are made manifest in your consciences.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>For we commend not ourselves again unto you, but give you occasion to glory on our behalf, that ye may have somewhat to answer them which glory in appearance, and not in heart.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>For whether we be beside ourselves, it is to God: or whether we be sober, it is for your cause.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">For the love of Christ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>constraineth</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> us; because we thus judge, that if one died for all, then were all dead:</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>And that he died for all, that they which live should not henceforth live unto themselves, but unto him which died for them, and rose again.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Wherefore henceforth know we no man after the flesh: yea, though we have known Christ after the flesh, yet now henceforth know we him no more.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Therefore</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> if any man be in Christ, he is a new creature: old things are passed away; behold, all things are become new.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">And all things are of God, who hath reconciled us to himself by Jesus </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Chrsit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, and hath given to us the ministry of reconciliation; </w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>To wit, that God was in Christ, reconciling the world unto himself, not imputing their trespasses unto them; and hath committed unto us the word of reconciliation.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Now then we are ambassadors for Christ, as though God did beseech you by us: we pray you in Christ’s stead, be ye reconciled to God.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>For he hath made him to be sin for us, who knew no sin; that we might be made the righteousness of God in him.</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

$d.Paragraphs.Item($prevIndex).Range.InsertParagraphAfter()
$curIndex = $prevIndex + 1
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>2 Corinthians Ch6</w:t></w:r></w:p>'
$d.Paragraphs.Item($curIndex).Range.InsertXML($frag)
$prevIndex = $curIndex

# --- drop the bootstrap paragraph now that list numId=$numId is minted ---
$d.Paragraphs.Item($bootstrapIndex).Range.Delete()

Write-Output ("Done. Minted numId=" + $numId + ", paragraphs now=" + $d.Paragraphs.Count)
